$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 5: the "this should be ignored" label that used to live at L5 (outside
#     the printed blueprint block) moves in next to the '#' in G5, landing on H5.
$ws.Range("L5").ClearContents()
$ws.Range("H5").Value = "this should be ignored"

# NOTE: allocate the two brand-new shared strings in the same order the
# target workbook uses (index 7 = "all # symbols...", index 8 = "x(3x1)"),
# by touching B17 before B8 below.
$ws.Range("B17").Value = "all # symbols should get ignored"

# --- Rows 8-12 get rebuilt as a taller "second layer" example (yellow s=2 block),
#     replacing the old rows 8-9. Stage the existing yellow style (s=2, still
#     present on A8) in a scratch cell first -- clearing A8:L12 below would
#     otherwise wipe the Office clipboard source before the paste runs -- then
#     clear the old content/formatting and re-apply that same style across the
#     whole new, larger block.
$ws.Range("A8").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A8:L12").Clear()

# Row 11 is brand new (the original file jumped from row 10 straight to row
# 12), so force its row-height metadata into existence now -- otherwise it
# won't pick up the sheet's explicit 17.25/customHeight row formatting when
# cells are written into it below.
$ws.Rows.Item(11).RowHeight = 17.25

$ws.Range("A30").Copy() | Out-Null
$ws.Range("A8:I12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("A30").Clear()

$ws.Range("B8").Value = "x(3x1)"
$ws.Range("I8").Value = "d"

$ws.Range("B11").Value = "x"
$ws.Range("C11").Value = "x"
$ws.Range("D11").Value = "x"
$ws.Range("I11").Value = "d"

$ws.Range("B12").Value = "x"
$ws.Range("C12").Value = "x"
$ws.Range("D12").Value = "x"
$ws.Range("I12").Value = "d"

# --- Rows 13, 15, 16, 17: former rows 10 and 12 shifted down, plus two brand
#     new rows (16 repeats the boundary row, 17 is the new "ignore #" note).
#     All four are brand new rows, so seed their row-height metadata first.
$ws.Rows.Item(13).RowHeight = 17.25
$ws.Rows.Item(15).RowHeight = 17.25
$ws.Rows.Item(16).RowHeight = 17.25
$ws.Rows.Item(17).RowHeight = 17.25

$ws.Range("A13").Value = "#"
$ws.Range("B13").Value = "#"
$ws.Range("C13").Value = "#"
$ws.Range("D13").Value = "#"
$ws.Range("E13").Value = "#"
$ws.Range("F13").Value = "#"
$ws.Range("G13").Value = "this should be ignored"

$ws.Range("B15").Value = "this should be ignored"

$ws.Range("A16").Value = "#"
$ws.Range("B16").Value = "#"
$ws.Range("C16").Value = "#"
$ws.Range("D16").Value = "#"
$ws.Range("E16").Value = "#"
$ws.Range("F16").Value = "#"
$ws.Range("G16").Value = "this should be ignored"

# --- Selection moves to B8 (matches the saved view state in the target file).
$ws.Range("B8").Select()
